# Applies the "more manual calculations and documentation" edit to the
# "Compiled Statement" sheet of the AAPL compiled income statement workbook.
#
# Summary of the change:
#   - Row 37 ("services"), Row 38 ("ebita ..."), Row 39 ("depreciation and
#     amortization") get corrected/re-entered values across D:AC (the data
#     that used to live in AD:AO is dropped - the trailing columns beyond
#     Q3 2022 / column AC were stale and are cleared out).
#   - Row 41 ("cash paid for income taxes, net") is a pure formula
#     (=col20+SUM(col37:col39)) so it recalculates on its own once rows
#     37-39 change; we just need to clear its stale AD:AO cells too.
#   - The sheet's used range / dimension naturally shrinks back down to
#     column AC once the AD:AO cells are cleared.
#   - The view is scrolled back to the top (no frozen/odd top-left cell)
#     and the on-save selection moves to roughly where the author was
#     last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compiled Statement")

# ---- Row 37 ("services") : D37:AC37 -------------------------------------
$row37 = @(2354,2484,10157,2745,2739,2665,2754,10903,3395,3040,2933,3179,12547,2816,2786,2752,2702,11056,2666,2797,2832,2989,11284,2697,2737,2805)
for ($i = 0; $i -lt $row37.Length; $i++) {
    $ws.Cells.Item(37, 4 + $i).Value = $row37[$i]
}

# ---- Row 38 ("ebita (income before taxes, interest, depreciation ...")  -
$row38 = @(2874,1839,11591,3551,2789,2479,1598,10417,4916,4581,2298,3468,15263,4393,3112,905,1091,9501,1787,8489,8260,6849,25385,5235,4066,2950)
for ($i = 0; $i -lt $row38.Length; $i++) {
    $ws.Cells.Item(38, 4 + $i).Value = $row38[$i]
}

# ---- Row 39 ("cash paid for income taxes, net") --------------------------
$row39 = @(449,636,2092,623,733,764,902,3022,836,926,801,860,3423,771,918,586,727,3002,619,708,543,817,2687,531,875,504)
for ($i = 0; $i -lt $row39.Length; $i++) {
    $ws.Cells.Item(39, 4 + $i).Value = $row39[$i]
}

# ---- Drop the stale trailing data beyond column AC (Q3 2022) -------------
# Rows 37-39 held manually-entered numbers out to column AO; row 41 carried
# the matching (now orphaned) formulas. Clearing these lets the sheet's
# used range / dimension and each row's "spans" collapse back to A:AC.
$ws.Range("AD37:AO39").ClearContents()
$ws.Range("AD41:AO41").ClearContents()

# Recalculate so row 41 (=col20+SUM(col37:col39)) reflects the new inputs.
$excel.CalculateFull()

# ---- View bookkeeping: match where the author left the cursor/selection --
$ws.Range("AK33:AX43").Select()
